# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the worker/period detail table (rows 16-31)
# Columns: C=Doc Trabajador, D=Nombre Trabajador, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @{ Row = 16; C = "9074280";  D = "LUIS ALBERTO ESTUPIÑAN GUIZA"; E = "1607"; F = 120000; G = 3961000 },
    @{ Row = 17; C = "9074280";  D = "LUIS ALBERTO ESTUPIÑAN GUIZA"; E = "1608"; F = 120000; G = 3961000 },
    @{ Row = 18; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1608"; F = 60000;  G = 5200000 },
    @{ Row = 19; C = "9074280";  D = "LUIS ALBERTO ESTUPIÑAN GUIZA"; E = "1609"; F = 120000; G = 3961000 },
    @{ Row = 20; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1609"; F = 60000;  G = 5200000 },
    @{ Row = 21; C = "22799499"; D = "MARIA CLAUDIA MONTIEL HERAZO"; E = "1610"; F = 128000; G = 16664000 },
    @{ Row = 22; C = "9074280";  D = "LUIS ALBERTO ESTUPIÑAN GUIZA"; E = "1610"; F = 120000; G = 3961000 },
    @{ Row = 23; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1610"; F = 60000;  G = 5200000 },
    @{ Row = 24; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1612"; F = 60000;  G = 5200000 },
    @{ Row = 25; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1701"; F = 60000;  G = 5200000 },
    @{ Row = 26; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1702"; F = 60000;  G = 5200000 },
    @{ Row = 27; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1703"; F = 60000;  G = 5200000 },
    @{ Row = 28; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1705"; F = 60000;  G = 5200000 },
    @{ Row = 29; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1706"; F = 183847; G = 5200000 },
    @{ Row = 30; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1707"; F = 183847; G = 5200000 },
    @{ Row = 31; C = "71753222"; D = "MAURICIO LOPEZ LONDOÑO";       E = "1708"; F = 183847; G = 5200000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
